$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week dates) ---
$ws.Range("A8").Value = "Volume 29   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/28/2022  Through  12/4/2022"

# --- Row 22: convert D22/E22 from "N/A" text markers to real numbers ---
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'

# --- Numeric data updates (weekly crime-stat refresh) ---
$ws.Range("L15").Value = -11.764705882352
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 23.076923076923
$ws.Range("I16").Value = 169
$ws.Range("J16").Value = 126
$ws.Range("K16").Value = 34.126984126984
$ws.Range("L16").Value = 10.457516339869
$ws.Range("M16").Value = -51.714285714285
$ws.Range("N16").Value = -85.963455149501
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -5.555555555555
$ws.Range("I17").Value = 259
$ws.Range("J17").Value = 204
$ws.Range("K17").Value = 26.960784313725
$ws.Range("L17").Value = 19.907407407407
$ws.Range("M17").Value = 36.315789473684
$ws.Range("N17").Value = -53.41726618705
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 23
$ws.Range("H18").Value = 15
$ws.Range("I18").Value = 280
$ws.Range("J18").Value = 213
$ws.Range("K18").Value = 31.455399061032
$ws.Range("L18").Value = -18.840579710144
$ws.Range("M18").Value = -36.936936936936
$ws.Range("N18").Value = -77.617905675459
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -57.894736842105
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -19.148936170212
$ws.Range("I19").Value = 622
$ws.Range("J19").Value = 479
$ws.Range("K19").Value = 29.853862212943
$ws.Range("L19").Value = 24.649298597194
$ws.Range("M19").Value = 33.190578158458
$ws.Range("N19").Value = 20.07722007722
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 9.090909090909
$ws.Range("I20").Value = 169
$ws.Range("J20").Value = 143
$ws.Range("K20").Value = 18.181818181818
$ws.Range("L20").Value = 36.290322580645
$ws.Range("M20").Value = 14.189189189189
$ws.Range("N20").Value = -80.72976054732
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -34.285714285714
$ws.Range("F21").Value = 106
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = -2.752293577981
$ws.Range("I21").Value = 1514
$ws.Range("J21").Value = 1186
$ws.Range("K21").Value = 27.655986509274
$ws.Range("L21").Value = 11.569638909358
$ws.Range("M21").Value = -5.904288377874
$ws.Range("N21").Value = -66.053811659192
$ws.Range("C22").Value = 2
$ws.Range("I22").Value = 39
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = 85.714285714285
$ws.Range("L22").Value = 77.272727272727
$ws.Range("M22").Value = -4.878048780487
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -80
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 147
$ws.Range("J23").Value = 159
$ws.Range("K23").Value = -7.54716981132
$ws.Range("L23").Value = -11.44578313253
$ws.Range("M23").Value = 16.666666666666
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -30.76923076923
$ws.Range("F24").Value = 64
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = -40.74074074074
$ws.Range("I24").Value = 1108
$ws.Range("J24").Value = 935
$ws.Range("K24").Value = 18.502673796791
$ws.Range("L24").Value = 24.076147816349
$ws.Range("M24").Value = -4.974271012006
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 83.333333333333
$ws.Range("I25").Value = 470
$ws.Range("J25").Value = 349
$ws.Range("K25").Value = 34.670487106017
$ws.Range("L25").Value = 57.718120805369
$ws.Range("M25").Value = 1.075268817204
$ws.Range("L26").Value = -4.545454545454
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 45
$ws.Range("J27").Value = 66
$ws.Range("K27").Value = -31.818181818181
$ws.Range("L27").Value = -10
$ws.Range("N28").Value = -87.951807228915
$ws.Range("N29").Value = -87.012987012987
$ws.Range("L30").Value = 211.111111111111
